$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRED Graph")

# New data rows appended to the bottom of the series (dates as Excel
# serial numbers, matching the existing column A date format; values in
# column B with the existing 2-decimal number format).
$newRows = @(
    @{ Row = 109; Date = 45786; Value = 2.08 },
    @{ Row = 110; Date = 45803; Value = 2.11 },
    @{ Row = 111; Date = 45807; Value = 2.11 },
    @{ Row = 112; Date = 45811; Value = 2.14 },
    @{ Row = 113; Date = 45820; Value = 2.11 },
    @{ Row = 114; Date = 45825; Value = 2.08 },
    @{ Row = 115; Date = 45838; Value = 1.95 }
)

foreach ($r in $newRows) {
    $aCell = $ws.Cells.Item($r.Row, 1)
    $aCell.Value = $r.Date
    $aCell.NumberFormat = "yyyy\-mm\-dd"

    $bCell = $ws.Cells.Item($r.Row, 2)
    $bCell.Value = $r.Value
    $bCell.NumberFormat = "0.00"
}

# Move the active selection to B116 (one row below the new last data row),
# matching the post-edit cursor position recorded in the sheet view.
$ws.Range("B116").Select()
